# Update "Escolher componente acessório" use case sheet:
#  - Normal scenario gains two new system-response steps
#    ("3. Apresenta preço do componente" and "6. Apresenta preço final"),
#    shifting the former steps 6/7 to 7/8.
#  - The exception note at the bottom now references the new step number (7 instead of 6).
#  - A stray "s" marker cell is present at G10 in the authored sheet.
#  - Selection/scroll position is updated to reflect the edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the previously-blank row with the new step 3 ---------------
$ws.Range("D9").Value = "3. Apresenta preço do componente"

# --- 2. Insert a new row for step "6. Apresenta preço final" ---------------
# Inserting before the old row 12 ("6. Confirma componte") pushes the
# remainder of the merged B6:B14 block down by one row (-> B6:B15) and
# shifts every block below it, exactly like doing it by hand in Excel.
$ws.Rows.Item(12).Insert()

# Match the formatting (borders/fill/font/alignment) of the surrounding
# rows, which all already share identical per-column styling.
$ws.Range("B11:D11").Copy()
$ws.Range("B12:D12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(12).RowHeight = 19
$ws.Range("D12").Value = "6. Apresenta preço final"

# --- 3. Renumber the steps that used to be 6/7 -> 7/8 -----------------------
$ws.Range("C13").Value = "7. Confirma componte"
$ws.Range("D14").Value = "8.  Adiciona componente"

# --- 4. Update the exception footnote: step 6 -> step 7 --------------------
$ws.Range("B23").Value = " Excepção 3 [Cliente não aceita componente] Passos 4.2, 5.2 e 7"

# --- 5. Stray helper cell left in the sheet by the author -------------------
$ws.Range("G10").Value = "s"

# --- 6. Update the view: scroll position + active selection ----------------
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G23").Select()
